$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 127.333336
$ws.Range("I9").Value = 132.8
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 132.8
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 36.19999999999999
$ws.Range("N9").Value = -438

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2691.2
$ws.Range("J17").Value = 2691.2
$ws.Range("L17").Value = 8073.599999999999
$ws.Range("N17").Value = -8409.599999999999

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9165.666999999999
$ws.Range("I74").Value = 8999
$ws.Range("J74").Value = 9499
$ws.Range("K74").Value = 8999
$ws.Range("L74").Value = 9499
$ws.Range("M74").Value = -8063
$ws.Range("N74").Value = -11371

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 9165.666999999999
$ws.Range("I77").Value = 8999
$ws.Range("J77").Value = 9499
$ws.Range("K77").Value = 44995
$ws.Range("L77").Value = 47495
$ws.Range("M77").Value = -40315
$ws.Range("N77").Value = -56855

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1749.75
$ws.Range("I112").Value = 1500
$ws.Range("K112").Value = 4500
$ws.Range("M112").Value = -3392

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 95000
$ws.Range("J140").Value = 95000
$ws.Range("L140").Value = 95000
$ws.Range("N140").Value = -105360

# ARM row 22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 17441.5
$ws.Range("J22").Value = 9750
$ws.Range("L22").Value = 9750
$ws.Range("N22").Value = -10348

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 11614.35
$ws.Range("J44").Value = 11614.35
$ws.Range("L44").Value = 11614.35
$ws.Range("N44").Value = -12590.35

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2318.5833
$ws.Range("I61").Value = 1601.25
$ws.Range("J61").Value = 3753.25
$ws.Range("K61").Value = 1601.25
$ws.Range("L61").Value = 3753.25
$ws.Range("M61").Value = -1389.25
$ws.Range("N61").Value = -4177.25

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2653.6428
$ws.Range("I74").Value = 2165.4614
$ws.Range("K74").Value = 2165.4614
$ws.Range("M74").Value = -1291.4614

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2653.6428
$ws.Range("I77").Value = 2165.4614
$ws.Range("K77").Value = 10827.307
$ws.Range("M77").Value = -6459.307000000001

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2131
$ws.Range("I97").Value = 1619.625
$ws.Range("K97").Value = 1619.625
$ws.Range("M97").Value = -1123.625

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2548.75
$ws.Range("I132").Value = 2398.3333
$ws.Range("K132").Value = 7194.999899999999
$ws.Range("M132").Value = -4664.999899999999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2318.5833
$ws.Range("I136").Value = 1601.25
$ws.Range("J136").Value = 3753.25
$ws.Range("K136").Value = 4803.75
$ws.Range("L136").Value = 11259.75
$ws.Range("M136").Value = -2253.75
$ws.Range("N136").Value = -16359.75

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4151
$ws.Range("I107").Value = 1444.3077
$ws.Range("J107").Value = 7083.25
$ws.Range("K107").Value = 1444.3077
$ws.Range("L107").Value = 7083.25
$ws.Range("M107").Value = 475.6922999999999
$ws.Range("N107").Value = -10923.25

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2173.111
$ws.Range("I134").Value = 1194.75
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 3584.25
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -1049.25
$ws.Range("N134").Value = -35070

# CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1251
$ws.Range("I3").Value = 875
$ws.Range("J3").Value = 1439
$ws.Range("K3").Value = 875
$ws.Range("L3").Value = 1439
$ws.Range("M3").Value = -762
$ws.Range("N3").Value = -1665

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 882285.3
$ws.Range("J41").Value = 1049730.6
$ws.Range("L41").Value = 1049730.6
$ws.Range("N41").Value = -1050586.6

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5816.143
$ws.Range("I58").Value = 1583
$ws.Range("K58").Value = 1583
$ws.Range("M58").Value = -1380

# CRP row 108
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 67541.664
$ws.Range("J108").Value = 90000
$ws.Range("L108").Value = 90000
$ws.Range("N108").Value = -97680

# CRP row 111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 61250
$ws.Range("J111").Value = 61250
$ws.Range("L111").Value = 61250
$ws.Range("N111").Value = -69430

# CRP row 112
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 37351
$ws.Range("J112").Value = 37351
$ws.Range("L112").Value = 37351
$ws.Range("N112").Value = -40305

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3031.762
$ws.Range("I132").Value = 2440.1667
$ws.Range("J132").Value = 6581.3335
$ws.Range("K132").Value = 7320.500100000001
$ws.Range("L132").Value = 19744.0005
$ws.Range("M132").Value = -4790.500100000001
$ws.Range("N132").Value = -24804.0005

# CRP row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 51878.855
$ws.Range("J133").Value = 51878.855
$ws.Range("L133").Value = 51878.855
$ws.Range("N133").Value = -56938.855

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5816.143
$ws.Range("I136").Value = 1583
$ws.Range("K136").Value = 4749
$ws.Range("M136").Value = -2199

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 61.666668
$ws.Range("I6").Value = 45
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 135
$ws.Range("L6").Value = 210
$ws.Range("M6").Value = -22
$ws.Range("N6").Value = -436

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 424.4762
$ws.Range("I17").Value = 79.57143000000001
$ws.Range("K17").Value = 238.71429
$ws.Range("M17").Value = -69.71429000000001

# CUL row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 558325.8
$ws.Range("I128").Value = 558325.8
$ws.Range("K128").Value = 1674977.4
$ws.Range("M128").Value = -1669997.4

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1415.75
$ws.Range("I132").Value = 757.6
$ws.Range("J132").Value = 1885.8572
$ws.Range("K132").Value = 6818.400000000001
$ws.Range("L132").Value = 16972.7148
$ws.Range("M132").Value = -4288.400000000001
$ws.Range("N132").Value = -22032.7148

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1200
$ws.Range("I97").Value = 1200
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -704
$ws.Range("N97").Value = -2192

# LTW row 33
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 2509.5
$ws.Range("I33").Value = 2509.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2509.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -2219.5
$ws.Range("N33").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2401.25
$ws.Range("I132").Value = 1866.6666
$ws.Range("J132").Value = 4005
$ws.Range("K132").Value = 5599.9998
$ws.Range("L132").Value = 12015
$ws.Range("M132").Value = -3069.9998
$ws.Range("N132").Value = -17075

# WVR row 10
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -331
$ws.Range("N10").ClearContents()

# WVR row 121
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# WVR row 127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 75000
$ws.Range("I127").Value = 75000
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 75000
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -70040
$ws.Range("N127").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7452.1816
$ws.Range("I136").Value = 7442.1113
$ws.Range("J136").Value = 7497.5
$ws.Range("K136").Value = 22326.3339
$ws.Range("L136").Value = 22492.5
$ws.Range("M136").Value = -19776.3339
$ws.Range("N136").Value = -27592.5
